$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 427.5
$ws.Range("I19").Value = 240
$ws.Range("J19").Value = 490
$ws.Range("K19").Value = 240
$ws.Range("L19").Value = 490
$ws.Range("M19").Value = -65
$ws.Range("N19").Value = -840
$ws.Range("H45").Value = 4583.4
$ws.Range("I45").Value = 2458.5
$ws.Range("J45").Value = 6000
$ws.Range("K45").Value = 7375.5
$ws.Range("L45").Value = 18000
$ws.Range("M45").Value = -7183.5
$ws.Range("N45").Value = -18384
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19436.28
$ws.Range("I32").Value = 5370.517
$ws.Range("K32").Value = 5370.517
$ws.Range("M32").Value = -5083.517
$ws.Range("H44").Value = 28049
$ws.Range("J44").Value = 28049
$ws.Range("L44").Value = 28049
$ws.Range("N44").Value = -29025
$ws.Range("H55").Value = 37926.5
$ws.Range("J55").Value = 37926.5
$ws.Range("L55").Value = 37926.5
$ws.Range("N55").Value = -38556.5
$ws.Range("H61").Value = 1227.7
$ws.Range("I61").Value = 1114.1111
$ws.Range("J61").Value = 2250
$ws.Range("K61").Value = 1114.1111
$ws.Range("L61").Value = 2250
$ws.Range("M61").Value = -902.1111000000001
$ws.Range("N61").Value = -2674
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H80").Value = 37105
$ws.Range("J80").Value = 37105
$ws.Range("L80").Value = 37105
$ws.Range("N80").Value = -39101
$ws.Range("H83").Value = 37105
$ws.Range("J83").Value = 37105
$ws.Range("L83").Value = 111315
$ws.Range("N83").Value = -121299
$ws.Range("H97").Value = 1486
$ws.Range("I97").Value = 930
$ws.Range("J97").Value = 1856.6666
$ws.Range("K97").Value = 930
$ws.Range("L97").Value = 1856.6666
$ws.Range("M97").Value = -434
$ws.Range("N97").Value = -2848.6666
$ws.Range("H136").Value = 1227.7
$ws.Range("I136").Value = 1114.1111
$ws.Range("J136").Value = 2250
$ws.Range("K136").Value = 3342.3333
$ws.Range("L136").Value = 6750
$ws.Range("M136").Value = -792.3333000000002
$ws.Range("N136").Value = -11850
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 26949.6
$ws.Range("J35").Value = 26949.6
$ws.Range("L35").Value = 26949.6
$ws.Range("N35").Value = -27569.6
$ws.Range("H80").Value = 304.14285
$ws.Range("I80").Value = 69.85714
$ws.Range("J80").Value = 421.2857
$ws.Range("K80").Value = 69.85714
$ws.Range("L80").Value = 421.2857
$ws.Range("M80").Value = 928.14286
$ws.Range("N80").Value = -2417.2857
$ws.Range("H82").Value = 23079.059
$ws.Range("J82").Value = 30008.455
$ws.Range("L82").Value = 30008.455
$ws.Range("N82").Value = -30774.455
$ws.Range("H83").Value = 304.14285
$ws.Range("I83").Value = 69.85714
$ws.Range("J83").Value = 421.2857
$ws.Range("K83").Value = 349.2857
$ws.Range("L83").Value = 2106.4285
$ws.Range("M83").Value = 4642.7143
$ws.Range("N83").Value = -12090.4285
$ws.Range("H85").Value = 23079.059
$ws.Range("J85").Value = 30008.455
$ws.Range("L85").Value = 30008.455
$ws.Range("N85").Value = -32660.455
$ws.Range("H94").Value = 758.56525
$ws.Range("I94").Value = 659.3158
$ws.Range("J94").Value = 1230
$ws.Range("K94").Value = 659.3158
$ws.Range("L94").Value = 1230
$ws.Range("M94").Value = -208.3158
$ws.Range("N94").Value = -2132
$ws.Range("H105").Value = 2980.4333
$ws.Range("I105").Value = 1861.2222
$ws.Range("K105").Value = 1861.2222
$ws.Range("M105").Value = -114.2221999999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 9128.200000000001
$ws.Range("J51").Value = 9128.200000000001
$ws.Range("L51").Value = 9128.200000000001
$ws.Range("N51").Value = -10600.2
$ws.Range("H61").Value = 9128.200000000001
$ws.Range("J61").Value = 9128.200000000001
$ws.Range("L61").Value = 9128.200000000001
$ws.Range("N61").Value = -9824.200000000001
$ws.Range("H76").Value = 3040
$ws.Range("I76").Value = 3040
$ws.Range("K76").Value = 3040
$ws.Range("M76").Value = -2725
$ws.Range("H79").Value = 3040
$ws.Range("I79").Value = 3040
$ws.Range("K79").Value = 3040
$ws.Range("M79").Value = -1948
$ws.Range("H132").Value = 2578.8
$ws.Range("I132").Value = 3227.4
$ws.Range("J132").Value = 1281.6
$ws.Range("K132").Value = 9682.200000000001
$ws.Range("L132").Value = 3844.8
$ws.Range("M132").Value = -7152.200000000001
$ws.Range("N132").Value = -8904.799999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2112280.2
$ws.Range("I2").Value = 111141
$ws.Range("J2").Value = 5714331
$ws.Range("K2").Value = 666846
$ws.Range("L2").Value = 34285986
$ws.Range("M2").Value = -666733
$ws.Range("N2").Value = -34286212
$ws.Range("H58").Value = 3464
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 3464
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 10392
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -10648
$ws.Range("H113").Value = 561
$ws.Range("I113").Value = 600
$ws.Range("J113").Value = 554.5
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 1663.5
$ws.Range("M113").Value = 370
$ws.Range("N113").Value = -6003.5
$ws.Range("H131").Value = 947.0816
$ws.Range("I131").Value = 586.6667
$ws.Range("J131").Value = 958.46313
$ws.Range("K131").Value = 1760.0001
$ws.Range("L131").Value = 2875.38939
$ws.Range("M131").Value = 3279.9999
$ws.Range("N131").Value = -12955.38939
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 966.6667
$ws.Range("I13").Value = 1000
$ws.Range("J13").Value = 950
$ws.Range("K13").Value = 1000
$ws.Range("L13").Value = 950
$ws.Range("M13").Value = -861
$ws.Range("N13").Value = -1228
$ws.Range("H17").Value = 9155.444
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 9155.444
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 9155.444
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -9491.444
$ws.Range("H97").Value = 1713.0769
$ws.Range("I97").Value = 1945
$ws.Range("J97").Value = 1342
$ws.Range("K97").Value = 1945
$ws.Range("L97").Value = 1342
$ws.Range("M97").Value = -1449
$ws.Range("N97").Value = -2334
$ws.Range("H107").Value = 706.62964
$ws.Range("I107").Value = 493.46667
$ws.Range("J107").Value = 973.0833
$ws.Range("K107").Value = 493.46667
$ws.Range("L107").Value = 973.0833
$ws.Range("M107").Value = 1426.53333
$ws.Range("N107").Value = -4813.0833
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3707.6924
$ws.Range("I62").Value = 3125
$ws.Range("J62").Value = 3966.6667
$ws.Range("K62").Value = 3125
$ws.Range("L62").Value = 3966.6667
$ws.Range("M62").Value = -2501
$ws.Range("N62").Value = -5214.6667
$ws.Range("H65").Value = 3707.6924
$ws.Range("I65").Value = 3125
$ws.Range("J65").Value = 3966.6667
$ws.Range("K65").Value = 15625
$ws.Range("L65").Value = 19833.3335
$ws.Range("M65").Value = -12505
$ws.Range("N65").Value = -26073.3335
